# Working on the ability to dump individual files, as well as fixing a
# major bug, on file 72.
#
# Rows 12 and 16 each carry an Abstract (column D) that still has one or
# more lingering `id="ParN">` anchor markers (left behind by the HTML ->
# plain-text conversion of the source abstract) plus stray blank lines.
# Clean the abstract up one marker at a time; in lock-step, widen the
# separator used between authors in the companion Authors cell (column
# E) by one extra space per cleanup pass - one pass per anchor marker
# found in that row's abstract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Repeat-Space($count) {
    $out = ""
    $i = 0
    while ($i -lt $count) {
        $out = $out + " "
        $i = $i + 1
    }
    return $out
}

function Clean-Row($row) {
    $dCell = $ws.Cells.Item($row, 4)   # column D - Abstract
    $eCell = $ws.Cells.Item($row, 5)   # column E - Authors

    $abstract = $dCell.Value2
    $authors  = $eCell.Value2

    # discover the current ",<spaces>" separator width used between authors
    $sepLen = 1
    if ($authors -match ",( +)") {
        $sepLen = $matches[1].Length
    }
    $sep = "," + (Repeat-Space $sepLen)

    $n = 1
    while ($abstract.Contains('id="Par' + $n + '">')) {
        # collapse runs of blank lines left behind in the abstract text
        $abstract = $abstract -replace "\n\n+", "`n"
        # drop this pass's anchor marker
        $abstract = $abstract.Replace('id="Par' + $n + '">', '')

        # widen the comma-space separator between authors by one space
        $newSep = $sep + " "
        $authors = $authors.Replace($sep, $newSep)
        $sep = $newSep

        $dCell.Value2 = $abstract
        $eCell.Value2 = $authors

        $n = $n + 1
    }
}

Clean-Row 12
Clean-Row 16
